$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K", formerly "Strike#") is being regenerated with new values
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
